$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading-percent values for data rows 2-25 (case with 380 kV),
# columns B,C,D,F,G,I,K,L,N (E,H,J,M stay 0 and are left untouched).
$cols = @(2, 3, 4, 6, 7, 9, 11, 12, 14)  # B,C,D,F,G,I,K,L,N as 1-based column indices

$data = @(
    @(15.91083539101749, 4.825501082034571, 7.368005756929478, 36.84228511115724, 3.705667921514321, 30.52496718431971, 13.5133341099533, 10.84291897573359, 21.93578902854646),
    @(15.7032247638344, 4.616673925000294, 7.374857746462713, 36.76022303406105, 3.708495499361875, 30.54431279294691, 13.37152055136399, 10.8206473206199, 21.9884425199562),
    @(15.57874467056088, 4.482434811026328, 7.37915758595267, 36.71827005374546, 3.710322960582109, 30.56179560916966, 13.28726692550226, 10.80920367700985, 22.02267694077661),
    @(15.52883436086294, 4.426254621284887, 7.380933318061249, 36.70330327348967, 3.711090707261032, 30.57032765101166, 13.25368048771652, 10.80510444565587, 22.03710720656925),
    @(15.52059792224647, 4.416837944213944, 7.381229603135665, 36.70094687624045, 3.711219584977751, 30.57182935125658, 13.24814972892815, 10.80445792498023, 22.03953231278205),
    @(15.57806817459061, 4.481683071500993, 7.379181438640585, 36.71805957474887, 3.710333221284125, 30.56190497834791, 13.28681089054387, 10.80914610542467, 22.02286961051066),
    @(15.83866849966395, 4.754766872366276, 7.370349209233431, 36.81224558986195, 3.706623965146029, 30.53047346509699, 13.46387488394632, 10.83477860637983, 21.95354884271058),
    @(16.37060137433345, 5.241194298875358, 7.353755043226753, 37.06338755094157, 3.700071081295411, 30.51336390459514, 13.83172176048583, 10.90259190215109, 21.83270561343635),
    @(16.7702756583993, 5.567195787929735, 7.34199185023995, 37.28764749439596, 3.695691149727022, 30.52799465940569, 14.11212841391229, 10.96287389446143, 21.75309521732725),
    @(16.9531999278719, 5.708464118466919, 7.336730519428955, 37.39809580219646, 3.693791872807915, 30.54055878670276, 14.24136828801561, 10.99250839522139, 21.718864421789),
    @(17.02256061460017, 5.760933091332507, 7.334750881200995, 37.44111167471556, 3.693085982748635, 30.54616523040112, 14.29050538250207, 11.00404240835529, 21.70618710874678),
    @(17.00761956621743, 5.749678815033429, 7.335176669447523, 37.4317947978379, 3.693237417343603, 30.54492005822773, 14.27991479324662, 11.0015445730081, 21.70890471598035),
    @(16.95890483170813, 5.712801417488359, 7.336567399745856, 37.40161099107301, 3.693733532177534, 30.54100303081721, 14.24540714001559, 10.99345108698173, 21.71781573795065),
    @(16.92907552488808, 5.690078892966301, 7.33742091196337, 37.38327705279079, 3.694039149954071, 30.53871422391989, 14.22429449762075, 10.98853404837268, 21.72331112499155),
    @(16.75833758993465, 5.557820718087793, 7.34233748054012, 37.2805973526278, 3.695817140567101, 30.52729238955767, 14.10371226516331, 10.96098121675689, 21.75537219302163),
    @(16.65383068837195, 5.474871982635909, 7.34537649447203, 37.21975186896613, 3.69693169222297, 30.52179827346869, 14.03013753774123, 10.94464084421378, 21.77554871315646),
    @(16.59382962441473, 5.426500946398528, 7.347132919595276, 37.18555012791631, 3.697581527387977, 30.51919443345543, 13.98798026041256, 10.93545087828359, 21.78734050312883),
    @(16.57353504478662, 5.410010306894759, 7.347729074916357, 37.17410714968063, 3.697803059697918, 30.51840837122941, 13.97373555119638, 10.93237531384257, 21.79136508846334),
    @(16.664944892969, 5.483770559303347, 7.34505211162529, 37.22614684342953, 3.696812138700251, 30.52232556887072, 14.03795336126942, 10.94635875915971, 21.77338155766592),
    @(16.97321159227155, 5.723661162926684, 7.336158565092851, 37.41044454595818, 3.693587450175996, 30.5421305341206, 14.25553789184707, 10.99581992049062, 21.71519061751616),
    @(17.17518284667327, 5.874456478376243, 7.330420130478416, 37.53782768724805, 3.691557563608948, 30.56002015235226, 14.39886817964597, 11.02996154930488, 21.6788214463075),
    @(17.067363352109, 5.794526298931933, 7.333476134840692, 37.46921384310269, 3.692633873001104, 30.55002003330473, 14.32228185476972, 11.01157545860994, 21.69808034862897),
    @(16.65991990242327, 5.47974963816934, 7.345198736233498, 37.22325324681199, 3.69686616060619, 30.52208545027222, 14.0344193822573, 10.94558145289647, 21.77436073037927),
    @(16.22486179797116, 5.115024459628429, 7.35816794623375, 36.98840405407282, 3.701767149658834, 30.51321769014746, 13.73024865778045, 10.88239064755522, 21.86378400103986)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($r, $cols[$j]).Value = $rowValues[$j]
    }
}
